$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 123
$ws.Range("A3").Value = 456

$ws.Range("A5").Select()
